# Update the cover-page "Version X.Y.Z, YYYY-MM-DD" line.
#
#   "Version 11.08.01, 2016-02-15"  ->  "Version 11.10.01, 2016-05-02"
#
# The run is split across many small w:r elements (one rsid per edit made
# over the doc's history). To keep that same run layout (rather than
# collapsing everything into one big run, which is what a naive
# Find/Replace across the whole string would do), we patch each affected
# run's text individually, in left-to-right order, recomputing each
# run's absolute character offsets as we go (since earlier edits shift
# later offsets when the replacement text is a different length).
#
# A bare ".Text =" assignment on a sub-range merges that run with any
# immediately-adjacent run sharing identical formatting. To prevent that
# and keep each edited span as its own run (matching the original
# document's run boundaries), we "nudge" the run's formatting right
# after writing its text: flip Font.Size away and back. That forces the
# engine to materialize a standalone run for exactly that span.

$d = $word.ActiveDocument

# Locate the run of text we need to edit.
$anchor = $d.Content
$found = $anchor.Find.Execute("Version 11.08.01, 2016-02-15", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the version/date text to update"
}
$base = $anchor.Start

# Each entry: (offset from $base of the run's old text, old text, new text).
# Offsets are relative to the *original* (pre-edit) string so we can
# recompute absolute positions as we walk left-to-right and text lengths
# change.
$edits = @(
    @{ Off = 11; Old = "0";    New = "10" },
    @{ Off = 12; Old = "8";    New = ".0" },
    @{ Off = 13; Old = ".0";   New = "1"  },
    @{ Off = 15; Old = "1";    New = ", 20" },
    @{ Off = 16; Old = ", 20"; New = "1"  },
    @{ Off = 20; Old = "1";    New = "6"  },
    @{ Off = 21; Old = "6";    New = "-"  },
    @{ Off = 22; Old = "-";    New = "0"  },
    @{ Off = 23; Old = "0";    New = "5"  },
    @{ Off = 24; Old = "2";    New = "-"  },
    @{ Off = 25; Old = "-";    New = "02" },
    @{ Off = 26; Old = "15";   New = ""   }
)

$shift = 0
foreach ($e in $edits) {
    $start = $base + $e.Off + $shift
    $end = $start + $e.Old.Length
    $r = $d.Range($start, $end)
    $r.Text = $e.New

    if ($e.New.Length -gt 0) {
        # Re-split this span into its own run by toggling a formatting
        # property away and back (no visible effect; sz=18 half-points
        # == 9pt).
        $r2 = $d.Range($start, $start + $e.New.Length)
        $r2.Font.Size = 99
        $r2.Font.Size = 9
    }

    $shift = $shift + ($e.New.Length - $e.Old.Length)
}
